# Update NATMI Bdnf-Ntrk2 LR-pair results with new TPM-based recomputation.
# Adds a new "Resolving-Mac" target-cluster category (shared string) and
# extends/refreshes the results table (rows 2-9, cols A-T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Bdnf"
$ws.Range("C2").Value = "Ntrk2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.601971
$ws.Range("H2").Value = 1.805913
$ws.Range("I2").Value = 0.09594307528308157
$ws.Range("J2").Value = 0.09594307528308157
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8630909999999999
$ws.Range("N2").Value = 2.589273
$ws.Range("O2").Value = 0.029079428547613
$ws.Range("P2").Value = 0.029079428547613
$ws.Range("Q2").Value = 0.519555752361
$ws.Range("R2").Value = 4.676001771249
$ws.Range("S2").Value = 0.002789969802332625
$ws.Range("T2").Value = 0.002789969802332625

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Bdnf"
$ws.Range("C3").Value = "Ntrk2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.601971
$ws.Range("H3").Value = 1.805913
$ws.Range("I3").Value = 0.09594307528308157
$ws.Range("J3").Value = 0.09594307528308157
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 24.72809866666667
$ws.Range("N3").Value = 74.184296
$ws.Range("O3").Value = 0.833143872773158
$ws.Range("P3").Value = 0.8331438727731579
$ws.Range("Q3").Value = 14.885598282472
$ws.Range("R3").Value = 133.970384542248
$ws.Range("S3").Value = 0.07993438530711323
$ws.Range("T3").Value = 0.07993438530711323

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bdnf"
$ws.Range("C4").Value = "Ntrk2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.601971
$ws.Range("H4").Value = 1.805913
$ws.Range("I4").Value = 0.09594307528308157
$ws.Range("J4").Value = 0.09594307528308157
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.080109666666666
$ws.Range("N4").Value = 12.240329
$ws.Range("O4").Value = 0.137467842346008
$ws.Range("P4").Value = 0.137467842346008
$ws.Range("Q4").Value = 2.456107696153
$ws.Range("R4").Value = 22.104969265377
$ws.Range("S4").Value = 0.01318908754720584
$ws.Range("T4").Value = 0.01318908754720584

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bdnf"
$ws.Range("C5").Value = "Ntrk2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.601971
$ws.Range("H5").Value = 1.805913
$ws.Range("I5").Value = 0.09594307528308157
$ws.Range("J5").Value = 0.09594307528308157
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.009167
$ws.Range("N5").Value = 0.027501
$ws.Range("O5").Value = 0.0003088563332209099
$ws.Range("P5").Value = 0.0003088563332209099
$ws.Range("Q5").Value = 0.005518268157
$ws.Range("R5").Value = 0.04966441341300001
$ws.Range("S5").Value = 0.00002963262642987029
$ws.Range("T5").Value = 0.00002963262642987029

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Bdnf"
$ws.Range("C6").Value = "Ntrk2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.672280666666667
$ws.Range("H6").Value = 17.016842
$ws.Range("I6").Value = 0.9040569247169185
$ws.Range("J6").Value = 0.9040569247169185
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8630909999999999
$ws.Range("N6").Value = 2.589273
$ws.Range("O6").Value = 0.029079428547613
$ws.Range("P6").Value = 0.029079428547613
$ws.Range("Q6").Value = 4.895694392874
$ws.Range("R6").Value = 44.061249535866
$ws.Range("S6").Value = 0.02628945874528038
$ws.Range("T6").Value = 0.02628945874528038

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Bdnf"
$ws.Range("C7").Value = "Ntrk2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.672280666666667
$ws.Range("H7").Value = 17.016842
$ws.Range("I7").Value = 0.9040569247169185
$ws.Range("J7").Value = 0.9040569247169185
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 24.72809866666667
$ws.Range("N7").Value = 74.184296
$ws.Range("O7").Value = 0.833143872773158
$ws.Range("P7").Value = 0.8331438727731579
$ws.Range("Q7").Value = 140.2647159903591
$ws.Range("R7").Value = 1262.382443913232
$ws.Range("S7").Value = 0.7532094874660449
$ws.Range("T7").Value = 0.7532094874660448

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Bdnf"
$ws.Range("C8").Value = "Ntrk2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.672280666666667
$ws.Range("H8").Value = 17.016842
$ws.Range("I8").Value = 0.9040569247169185
$ws.Range("J8").Value = 0.9040569247169185
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.080109666666666
$ws.Range("N8").Value = 12.240329
$ws.Range("O8").Value = 0.137467842346008
$ws.Range("P8").Value = 0.137467842346008
$ws.Range("Q8").Value = 23.14352718011311
$ws.Range("R8").Value = 208.291744621018
$ws.Range("S8").Value = 0.1242787547988022
$ws.Range("T8").Value = 0.1242787547988022

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Bdnf"
$ws.Range("C9").Value = "Ntrk2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.672280666666667
$ws.Range("H9").Value = 17.016842
$ws.Range("I9").Value = 0.9040569247169185
$ws.Range("J9").Value = 0.9040569247169185
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.009167
$ws.Range("N9").Value = 0.027501
$ws.Range("O9").Value = 0.0003088563332209099
$ws.Range("P9").Value = 0.0003088563332209099
$ws.Range("Q9").Value = 0.05199779687133333
$ws.Range("R9").Value = 0.467980171842
$ws.Range("S9").Value = 0.0002792237067910396
$ws.Range("T9").Value = 0.0002792237067910396
